# feat: add 2022-Q1 data
#
# Plan:
#  1. The existing "总计" (totals) sheet (currently the 5th sheet) is
#     repurposed to hold the new "2022-Q1" fund-holdings detail data
#     (same shape as the other quarterly sheets: 基金代码/基金名称/...).
#     Its B1:D1 header and A2:A5 row-index cells already carry the
#     workbook's "bold header / index column" style, so only content is
#     cleared (not formatting) there; the newly introduced cells
#     (E1:H1 header, A6:A9 index column) get that same style copied
#     across from an existing styled cell.
#  2. A brand-new "总计" sheet is appended after it, containing the
#     historical roll-up table (日期/持有数量(只)/持有市值(亿元)) with a
#     new leading row for 2022-Q1, styled to match the other sheets.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# A cell carrying the workbook's standard bold/bordered header-and-index
# style, used as the source for copying formatting onto newly added
# cells that have no pre-existing style to inherit.
$styleSource = $wb.Worksheets.Item(4).Range("B1")

# ---------------------------------------------------------------------
# Step 1: rename the current "总计" sheet to "2022-Q1" and replace its
# contents with the quarterly fund-holding detail rows.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(5)
$q1.Cells.ClearContents()
$q1.Name = "2022-Q1"

# New header cells E1:H1 and new index cells A6:A9 need the header/index
# style copied in (B1:D1 and A2:A5 already carry it from the old sheet).
$styleSource.Copy()
$q1.Range("E1:H1").PasteSpecial($xlPasteFormats)
$styleSource.Copy()
$q1.Range("A6:A9").PasteSpecial($xlPasteFormats)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# 基金代码 (fund code, has leading zeros) and 基金规模 / 股票总仓位 /
# 仓位占比 / 持有市值(亿元) are stored as text in this workbook's
# convention (not numbers) - force text formatting so values like
# "001838" / "32.13" are not reinterpreted as numeric.
$q1.Range("B2:B9").NumberFormat = "@"
$q1.Range("D2:G9").NumberFormat = "@"

$q1Data = @(
    @(0, "001838", "国投瑞银国家安全灵活配置混合", "32.13", "94.68", "5.59", "1.7961", 9),
    @(1, "001907", "国投瑞银境煊灵活配置混合A",   "2.61",  "90.44", "4.72", "0.1232", 4),
    @(2, "001908", "国投瑞银境煊灵活配置混合C",   "1.75",  "90.44", "4.72", "0.0826", 4),
    @(3, "257050", "国联安主题驱动混合",          "1.50",  "65.37", "2.89", "0.0434", 9),
    @(4, "015309", "国投瑞银境煊灵活配置混合E",   "0.33",  "90.44", "4.72", "0.0156", 4),
    @(5, "001899", "东海中证社会发展安全产业主题指数", "0.21", "90.30", "2.72", "0.0057", 5),
    @(6, "005104", "富荣福康混合A",                "0.08",  "87.88", "3.08", "0.0025", 5),
    @(7, "005105", "富荣福康混合C",                "0.04",  "87.88", "3.08", "0.0012", 5)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Range("A$r").Value = $row[0]
    $q1.Range("B$r").Value = $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = $row[3]
    $q1.Range("E$r").Value = $row[4]
    $q1.Range("F$r").Value = $row[5]
    $q1.Range("G$r").Value = $row[6]
    $q1.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: add the new "总计" sheet after "2022-Q1" with the updated
# roll-up table (2022-Q1 prepended).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$styleSource.Copy()
$total.Range("B1:D1").PasteSpecial($xlPasteFormats)
$styleSource.Copy()
$total.Range("A2:A6").PasteSpecial($xlPasteFormats)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @(0, "2022-Q1", 8, 2.07),
    @(1, "2021-Q4", 7, 4.88),
    @(2, "2021-Q3", 5, 0.23),
    @(3, "2021-Q2", 2, 0.2),
    @(4, "2020-Q4", 1, 0.57)
)

$r = 2
foreach ($row in $totalData) {
    $total.Range("A$r").Value = $row[0]
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = $row[2]
    $total.Range("D$r").Value = $row[3]
    $r = $r + 1
}

# Restore the originally active sheet/selection so this data-only edit
# doesn't also shuffle the workbook's view state.
$wb.Worksheets.Item(1).Activate() | Out-Null
$wb.Worksheets.Item(1).Range("A1").Select() | Out-Null
